# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) counts pulled from the regenerated source data (replacing
# the old Strike# counts that previously lived in column G).
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 5
    10 = 2
    11 = 2
    12 = 0
    13 = 0
    14 = 1
    15 = 3
    16 = 1
    17 = 2
    18 = 5
    19 = 2
    20 = 1
    21 = 2
    22 = 3
    23 = 1
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 2
    31 = 0
    32 = 6
    33 = 4
    34 = 2
    35 = 1
    36 = 1
    37 = 1
    38 = 3
    39 = 4
}

foreach ($row in $kValues.Keys | Sort-Object) {
    $ws.Range("G$row").Value = $kValues[$row]
}
